$d = $word.ActiveDocument

$d.Content.Find.Execute("618×8=", $true, $false, $false, $false, $false, $true, 1, $false, "922×5=", 2) | Out-Null
$d.Content.Find.Execute("588×4=", $true, $false, $false, $false, $false, $true, 1, $false, "470×2=", 2) | Out-Null
$d.Content.Find.Execute("917×9=", $true, $false, $false, $false, $false, $true, 1, $false, "104×3=", 2) | Out-Null
$d.Content.Find.Execute("884×8=", $true, $false, $false, $false, $false, $true, 1, $false, "724×9=", 2) | Out-Null
$d.Content.Find.Execute("799×6=", $true, $false, $false, $false, $false, $true, 1, $false, "761×7=", 2) | Out-Null
$d.Content.Find.Execute("601×8=", $true, $false, $false, $false, $false, $true, 1, $false, "615×8=", 2) | Out-Null
$d.Content.Find.Execute("855×6=", $true, $false, $false, $false, $false, $true, 1, $false, "380×2=", 2) | Out-Null
$d.Content.Find.Execute("827×7=", $true, $false, $false, $false, $false, $true, 1, $false, "568×3=", 2) | Out-Null
$d.Content.Find.Execute("354×5=", $true, $false, $false, $false, $false, $true, 1, $false, "176×5=", 2) | Out-Null
$d.Content.Find.Execute("886×8=", $true, $false, $false, $false, $false, $true, 1, $false, "339×5=", 2) | Out-Null
$d.Content.Find.Execute("559×7=", $true, $false, $false, $false, $false, $true, 1, $false, "258×2=", 2) | Out-Null
$d.Content.Find.Execute("451×4=", $true, $false, $false, $false, $false, $true, 1, $false, "726×7=", 2) | Out-Null
$d.Content.Find.Execute("647×5=", $true, $false, $false, $false, $false, $true, 1, $false, "205×7=", 2) | Out-Null
$d.Content.Find.Execute("330×7=", $true, $false, $false, $false, $false, $true, 1, $false, "589×5=", 2) | Out-Null
$d.Content.Find.Execute("983×3=", $true, $false, $false, $false, $false, $true, 1, $false, "959×3=", 2) | Out-Null
$d.Content.Find.Execute("328×7=", $true, $false, $false, $false, $false, $true, 1, $false, "854×5=", 2) | Out-Null
$d.Content.Find.Execute("674×8=", $true, $false, $false, $false, $false, $true, 1, $false, "189×8=", 2) | Out-Null
$d.Content.Find.Execute("456×9=", $true, $false, $false, $false, $false, $true, 1, $false, "645×7=", 2) | Out-Null
$d.Content.Find.Execute("145×4=", $true, $false, $false, $false, $false, $true, 1, $false, "786×7=", 2) | Out-Null
$d.Content.Find.Execute("604×3=", $true, $false, $false, $false, $false, $true, 1, $false, "395×5=", 2) | Out-Null
$d.Content.Find.Execute("188×2=", $true, $false, $false, $false, $false, $true, 1, $false, "426×4=", 2) | Out-Null
$d.Content.Find.Execute("773×9=", $true, $false, $false, $false, $false, $true, 1, $false, "419×4=", 2) | Out-Null
$d.Content.Find.Execute("218×3=", $true, $false, $false, $false, $false, $true, 1, $false, "484×8=", 2) | Out-Null
$d.Content.Find.Execute("881×3=", $true, $false, $false, $false, $false, $true, 1, $false, "372×4=", 2) | Out-Null
$d.Content.Find.Execute("564×2=", $true, $false, $false, $false, $false, $true, 1, $false, "487×4=", 2) | Out-Null
